$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 136; everything from 136 downward shifts to 137+.
$ws.Rows.Item(136).Insert()

# Copy the constant (unchanged) columns from the row that is now 137
# (formerly row 136) into the newly inserted row 136.
$ws.Range("A137:C137").Copy()
$ws.Range("A136").PasteSpecial()
$ws.Range("E137:I137").Copy()
$ws.Range("E136").PasteSpecial()
$ws.Range("O137:O137").Copy()
$ws.Range("O136").PasteSpecial()
$ws.Range("R137:R137").Copy()
$ws.Range("R136").PasteSpecial()

# Match style of the date cell (D column) from the row below.
$ws.Range("D137").Copy()
$ws.Range("D136").PasteSpecial()

# New record's own values.
$ws.Cells.Item(136, 4).Value = 44855
$ws.Cells.Item(136, 10).Value = 560
$ws.Cells.Item(136, 11).Value = 11500
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = 11750
$ws.Cells.Item(136, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(136, 16).Value = 294
$ws.Cells.Item(136, 17).Value = 40
